$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductBacklog")

# Swap values between row 12 (D12/F12) and row 17 (D17/F17)
$ws.Range("D12").Value = 9
$ws.Range("F12").Value = 2

$ws.Range("D17").Value = 15
$ws.Range("F17").Value = 3

# Update selected cell to H10
$ws.Range("H10").Select()
